$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("newaccountvalid")
$ws.Range("A2").Value = "emailemail@email.com"
$ws.Range("A3").Value = "emailemail@email.com"
